# Commit: "one more update to files for formatting"
#
# The sheet gained 4 new (blank/zeroed) data rows at the top of the data
# block (rows 2:5), pushing the existing 67 data rows down from rows
# 2:68 to rows 6:72. The newly inserted rows 2:5 are filled with 0 across
# every column (A:XFD, i.e. the full row width), matching the selection
# left behind by the edit (A2:XFD5).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Social_1002")

# Insert 4 new rows above the current row 2, shifting the existing data
# (previously rows 2:68) down to rows 6:72.
$ws.Rows("2:5").Insert()

# The newly inserted rows start out blank; zero-fill them across the full
# row width.
$ws.Range("A2:XFD5").Value = 0

# Leave the same range selected, matching the saved selection state.
$ws.Range("A2:XFD5").Select() | Out-Null
